# 036 : speedbal driver
#
# The "speedbal" (Speed Ball) entry is being promoted out of the "ALL"
# sheet and into the "Playable (untested)" sheet (it now has a working
# driver). Concretely:
#   1. Append the speedbal row as the new last row (308) of
#      "Playable (untested)".
#   2. Delete the old speedbal row (row 1601) from "ALL", which shifts
#      every following row up by one and shrinks the used range from
#      M1737 to M1736.
#   3. Re-point the AutoFilter on "ALL" and the two workbook-level
#      defined names (_FilterDatabase / LIST) at the new, smaller range.
#   4. Refresh the selections shown in each sheet's view.

$wb = $excel.ActiveWorkbook
$wsPlay = $wb.Worksheets.Item("Playable (untested)")
$wsAll  = $wb.Worksheets.Item("ALL")

# --- 1) add the speedbal row to "Playable (untested)" as row 308 -----------
$wsPlay.Range("A308").Value = 308
$wsPlay.Range("B308").Value = "speedbal"
$wsPlay.Range("C308").Value = "speedbal.c"
$wsPlay.Range("D308").Value = "Z80"
$wsPlay.Range("E308").Value = "Z80"
$wsPlay.Range("H308").Value = "1xYM-3812"
$wsPlay.Range("M308").Value = "Speed Ball"

# Update the view: selection now covers the freshly added row.
$wsPlay.Activate()
$wsPlay.Range("A306:A308").Select()

# --- 2) remove the now-duplicated speedbal row from "ALL" ------------------
$wsAll.Rows.Item(1601).Delete()

# --- 3) re-apply the AutoFilter over the new, one-row-smaller range --------
$wsAll.AutoFilterMode = $false
$wsAll.Range("A1:M1736").AutoFilter()

# Update the view: scroll position + selection after the deletion.
$wsAll.Activate()
$wsAll.Range("F1609").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1346

# --- 4) fix up the workbook-level defined names -----------------------------
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
  $n = $names.Item($i)
  if ($n.Name() -eq "ALL!_FilterDatabase") {
    $n.RefersTo = "=ALL!`$A`$1:`$M`$1736"
  }
  if ($n.Name() -eq "ALL!LIST") {
    $n.RefersTo = "=ALL!`$B`$1:`$M`$1736"
  }
}
